$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '70.283.26'
$ws.Range("E2").Value = '  -1.60%  '

# Row 3
$ws.Range("D3").Value = '3.798.23'
$ws.Range("E3").Value = '  +3.63%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '619.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.18%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.27%  '

# Row 7
$ws.Range("D7").Value = '3.797.50'
$ws.Range("E7").Value = '  +3.76%  '

# Row 8
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("E9").Value = '  +0.14%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.170'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.25%  '

# Row 11
$ws.Range("E11").Value = '  -3.55%  '

# Row 12
$ws.Range("E12").Value = '  -1.44%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.20'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.97%  '

# Row 14
$ws.Range("E14").Value = '  -0.27%  '

# Row 15
$ws.Range("D15").Value = '4.432.91'
$ws.Range("E15").Value = '  +3.61%  '

# Row 16
$ws.Range("D16").Value = '3.800.80'
$ws.Range("E16").Value = '  +3.43%  '

# Row 17
$ws.Range("D17").Value = '70.311.45'
$ws.Range("E17").Value = '  -1.50%  '

# Row 18
$ws.Range("E18").Value = '  -0.59%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.61'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.59%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '515.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.28%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.20%  '

# Row 22
$ws.Range("E22").Value = '  +3.58%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.730'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.97%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.54'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.59%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '88.26'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.45%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.33'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.94%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.83%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000138'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +22.41%  '

# Row 30
$ws.Range("E30").Value = '  -2.33%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.85'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.73%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.84'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.26%  '

# Row 33
$ws.Range("E33").Value = '  -1.58%  '

# Row 34
$ws.Range("E34").Value = '  -2.14%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.02%  '

# Row 36
$ws.Range("E36").Value = '  +1.05%  '

# Row 37
$ws.Range("E37").Value = '  +2.94%  '

# Row 38
$ws.Range("E38").Value = '  +0.76%  '

# Row 39
$ws.Range("E39").Value = '  +1.73%  '

# Row 40
$ws.Range("E40").Value = '  +3.40%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.16'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.35%  '

# Row 42
$ws.Range("B42").Value = 'Arweave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '44.35'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.65%  '

# Row 43
$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.79'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.50%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '423.74'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.22%  '

# Row 45
$ws.Range("D45").Value = '3.067.79'
$ws.Range("E45").Value = '  -2.82%  '

# Row 46
$ws.Range("E46").Value = '  -2.41%  '

# Row 47
$ws.Range("E47").Value = '  -0.40%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.67'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.45%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '135.83'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.71%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.49'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.37%  '
